$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Column D width: 30 -> 28 ---
$ws.Columns.Item(4).ColumnWidth = 27.17

# --- Remove hyperlinks belonging to rows 5-11 (keep F2/F3/F4 hyperlinks intact) ---
$toDelete = @()
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Row -gt 4) {
        $toDelete += $hl
    }
}
for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
    $toDelete[$i].Delete()
}

# --- Update the remaining hyperlink targets for rows 2-4 ---
foreach ($hl in $ws.Hyperlinks) {
    $r = $hl.Range.Row
    if ($r -eq 2) { $hl.Address = "https://www.lancers.jp/work/detail/5428507" }
    if ($r -eq 3) { $hl.Address = "https://www.lancers.jp/work/detail/5429882" }
    if ($r -eq 4) { $hl.Address = "https://www.lancers.jp/work/detail/5435519" }
}

# --- Delete rows 5-11 entirely (dimension shrinks to A1:H4) ---
$ws.Rows("5:11").Delete()

# --- Row 2 ---
$ws.Range("A2").Value = "2025-11-17 06:28:00"
$ws.Range("B2").Value = "【Next.js × TypeScript × Tailwind】コンポーネント制作パートナー募集!"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5428507"
$ws.Range("G2").Value = 528
$ws.Range("H2").Value = "🔥AI,Next.js"

# --- Row 3 ---
$ws.Range("A3").Value = "2025-11-17 06:28:00"
$ws.Range("B3").Value = "UTAGE構築代行|ヒアリングから構築までお任せしたいです。"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5429882"
$ws.Range("G3").Value = 18
$ws.Range("H3").Clear()

# --- Row 4 ---
$ws.Range("A4").Value = "2025-11-17 06:28:00"
$ws.Range("B4").Value = "ロリポップ!レンタルサーバーの不具合を解決したい"
$ws.Range("D4").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5435519"
$ws.Range("G4").Value = 10
$ws.Range("H4").Clear()
